# Auto-generated edit script: updates market price / profit columns (H-N)
# across the 8 crafting-job sheets to reflect refreshed marketboard data.
$wb = $excel.ActiveWorkbook

# ALC row 69
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 7291.75
$ws.Range("I69").Value = 6000
$ws.Range("J69").Value = 8583.5
$ws.Range("K69").Value = 18000
$ws.Range("L69").Value = 25750.5
$ws.Range("M69").Value = -17126
$ws.Range("N69").Value = -27498.5

# ALC row 72
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H72").Value = 7291.75
$ws.Range("I72").Value = 6000
$ws.Range("J72").Value = 8583.5
$ws.Range("K72").Value = 54000
$ws.Range("L72").Value = 77251.5
$ws.Range("M72").Value = -49632
$ws.Range("N72").Value = -85987.5

# ALC row 96
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 1129.7646
$ws.Range("I96").Value = 826.61536
$ws.Range("K96").Value = 2479.84608
$ws.Range("M96").Value = -1106.84608

# ALC row 100
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1353.2
$ws.Range("I100").Value = 1277.4117
$ws.Range("J100").Value = 1782.6666
$ws.Range("K100").Value = 1277.4117
$ws.Range("L100").Value = 1782.6666
$ws.Range("M100").Value = -736.4117000000001
$ws.Range("N100").Value = -2864.6666

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2418.75
$ws.Range("I137").Value = 1515.7693
$ws.Range("K137").Value = 4547.3079
$ws.Range("M137").Value = -1997.3079

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4561.212
$ws.Range("J138").Value = 3797.7083
$ws.Range("L138").Value = 11393.1249
$ws.Range("N138").Value = -21673.1249

# ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 3016
$ws.Range("I141").Value = 2024.75
$ws.Range("K141").Value = 6074.25
$ws.Range("M141").Value = -894.25

# ARM row 5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 486.4
$ws.Range("I5").Value = 147.33333
$ws.Range("K5").Value = 147.33333
$ws.Range("M5").Value = -35.33332999999999

# ARM row 97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 2484
$ws.Range("J97").Value = 1000
$ws.Range("L97").Value = 1000
$ws.Range("N97").Value = -1992

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2359.3635
$ws.Range("I132").Value = 2400.5
$ws.Range("J132").Value = 1948
$ws.Range("K132").Value = 7201.5
$ws.Range("L132").Value = 5844
$ws.Range("M132").Value = -4671.5
$ws.Range("N132").Value = -10904

# BSM row 4
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 486.4
$ws.Range("I4").Value = 147.33333
$ws.Range("K4").Value = 147.33333
$ws.Range("M4").Value = -32.33332999999999

# BSM row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 6636.815
$ws.Range("I20").Value = 6166.95
$ws.Range("J20").Value = 7979.2856
$ws.Range("K20").Value = 6166.95
$ws.Range("L20").Value = 7979.2856
$ws.Range("M20").Value = -5919.95
$ws.Range("N20").Value = -8473.285599999999

# BSM row 100
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H100").Value = 17693.6
$ws.Range("J100").Value = 17693.6
$ws.Range("L100").Value = 17693.6
$ws.Range("N100").Value = -19857.6

# BSM row 106
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H106").Value = 49835.5
$ws.Range("J106").Value = 49835.5
$ws.Range("L106").Value = 49835.5
$ws.Range("N106").Value = -52359.5

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4779.069
$ws.Range("I134").Value = 4196.269
$ws.Range("K134").Value = 12588.807
$ws.Range("M134").Value = -10053.807

# BSM row 137
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H137").Value = 75000
$ws.Range("J137").Value = 75000
$ws.Range("L137").Value = 75000
$ws.Range("N137").Value = -85200

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 76928296
$ws.Range("I31").Value = 125002856
$ws.Range("K31").Value = 125002856
$ws.Range("M31").Value = -125002561

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 76928296
$ws.Range("I34").Value = 125002856
$ws.Range("K34").Value = 125002856
$ws.Range("M34").Value = -125002654

# CRP row 36
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 1500
$ws.Range("I36").Value = 1000
$ws.Range("J36").Value = 1666.6666
$ws.Range("K36").Value = 1000
$ws.Range("L36").Value = 1666.6666
$ws.Range("M36").Value = -612
$ws.Range("N36").Value = -2442.6666

# CRP row 40
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H40").Value = 1500
$ws.Range("I40").Value = 1000
$ws.Range("J40").Value = 1666.6666
$ws.Range("K40").Value = 1000
$ws.Range("L40").Value = 1666.6666
$ws.Range("M40").Value = -840
$ws.Range("N40").Value = -1986.6666

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 7819.7393
$ws.Range("I58").Value = 3285.6
$ws.Range("K58").Value = 3285.6
$ws.Range("M58").Value = -3082.6

# CRP row 86
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 3488.2778
$ws.Range("I86").Value = 3196.923
$ws.Range("K86").Value = 3196.923
$ws.Range("M86").Value = -2073.923

# CRP row 89
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 3488.2778
$ws.Range("I89").Value = 3196.923
$ws.Range("K89").Value = 15984.615
$ws.Range("M89").Value = -10368.615

# CRP row 94
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 3396.625
$ws.Range("J94").Value = 4245.5
$ws.Range("L94").Value = 4245.5
$ws.Range("N94").Value = -5147.5

# CRP row 106
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").Value = $null

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4715.737
$ws.Range("I132").Value = 2789.375
$ws.Range("J132").Value = 14989.667
$ws.Range("K132").Value = 8368.125
$ws.Range("L132").Value = 44969.001
$ws.Range("M132").Value = -5838.125
$ws.Range("N132").Value = -50029.001

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2829.3704
$ws.Range("I134").Value = 1774.9474
$ws.Range("J134").Value = 5333.625
$ws.Range("K134").Value = 5324.8422
$ws.Range("L134").Value = 16000.875
$ws.Range("M134").Value = -2789.8422
$ws.Range("N134").Value = -21070.875

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 7819.7393
$ws.Range("I136").Value = 3285.6
$ws.Range("K136").Value = 9856.799999999999
$ws.Range("M136").Value = -7306.799999999999

# CUL row 2
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 33.5
$ws.Range("I2").Value = 38.307693
$ws.Range("J2").Value = 12.666667
$ws.Range("K2").Value = 229.846158
$ws.Range("L2").Value = 76.00000199999999
$ws.Range("M2").Value = -116.846158
$ws.Range("N2").Value = -302.000002

# CUL row 18
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 653.7
$ws.Range("I18").Value = 393
$ws.Range("K18").Value = 1179
$ws.Range("M18").Value = -1010

# CUL row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 71434000
$ws.Range("I68").Value = 125005000
$ws.Range("J68").Value = 5999.6665
$ws.Range("K68").Value = 375015000
$ws.Range("L68").Value = 17998.9995
$ws.Range("M68").Value = -375014189
$ws.Range("N68").Value = -19620.9995

# CUL row 71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 71434000
$ws.Range("I71").Value = 125005000
$ws.Range("J71").Value = 5999.6665
$ws.Range("K71").Value = 1125045000
$ws.Range("L71").Value = 53996.9985
$ws.Range("M71").Value = -1125040944
$ws.Range("N71").Value = -62108.9985

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6795.952
$ws.Range("I70").Value = 4122.4
$ws.Range("K70").Value = 4122.4
$ws.Range("M70").Value = -3852.4

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 6795.952
$ws.Range("I73").Value = 4122.4
$ws.Range("K73").Value = 4122.4
$ws.Range("M73").Value = -3186.4

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2717.4285
$ws.Range("I80").Value = 2362
$ws.Range("J80").Value = 2859.6
$ws.Range("K80").Value = 2362
$ws.Range("L80").Value = 2859.6
$ws.Range("M80").Value = -1364
$ws.Range("N80").Value = -4855.6

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2717.4285
$ws.Range("I83").Value = 2362
$ws.Range("J83").Value = 2859.6
$ws.Range("K83").Value = 11810
$ws.Range("L83").Value = 14298
$ws.Range("M83").Value = -6818
$ws.Range("N83").Value = -24282

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1696.5
$ws.Range("I122").Value = 1696.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5089.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2639.5
$ws.Range("N122").Value = $null

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2714.4285
$ws.Range("I132").Value = 2126.1667
$ws.Range("K132").Value = 6378.500100000001
$ws.Range("M132").Value = -3848.500100000001

# LTW row 55
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 259.94736
$ws.Range("I55").Value = 218.09091
$ws.Range("K55").Value = 218.09091
$ws.Range("M55").Value = -45.09091000000001

# LTW row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 5309.227
$ws.Range("I100").Value = 4254.5386
$ws.Range("K100").Value = 4254.5386
$ws.Range("M100").Value = -3713.5386

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2603.0908
$ws.Range("I122").Value = 2404.2222
$ws.Range("K122").Value = 7212.6666
$ws.Range("M122").Value = -4762.6666

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 11444.82
$ws.Range("I132").Value = 10730.263
$ws.Range("K132").Value = 32190.789
$ws.Range("M132").Value = -29660.789

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 7897.4
$ws.Range("I136").Value = 9865.666999999999
$ws.Range("J136").Value = 4945
$ws.Range("K136").Value = 29597.001
$ws.Range("L136").Value = 14835
$ws.Range("M136").Value = -27047.001
$ws.Range("N136").Value = -19935

# WVR row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2995.5715
$ws.Range("J81").Value = 7117.3335
$ws.Range("L81").Value = 14234.667
$ws.Range("N81").Value = -16356.667

# WVR row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 2995.5715
$ws.Range("J84").Value = 7117.3335
$ws.Range("L84").Value = 71173.33499999999
$ws.Range("N84").Value = -81781.33499999999

# WVR row 96
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2895.4546
$ws.Range("I96").Value = 2724
$ws.Range("J96").Value = 3195.5
$ws.Range("K96").Value = 2724
$ws.Range("L96").Value = 3195.5
$ws.Range("M96").Value = -1351
$ws.Range("N96").Value = -5941.5

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3873.682
$ws.Range("I126").Value = 3787.8157
$ws.Range("K126").Value = 11363.4471
$ws.Range("M126").Value = -8893.447100000001

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I132").Value = 3144.4517
$ws.Range("J132").Value = 3972.818
$ws.Range("K132").Value = 9433.355100000001
$ws.Range("L132").Value = 11918.454
$ws.Range("M132").Value = -6903.355100000001
$ws.Range("N132").Value = -16978.454

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3491.0408
$ws.Range("I136").Value = 2242.2163
$ws.Range("J136").Value = 7341.5835
$ws.Range("K136").Value = 6726.6489
$ws.Range("L136").Value = 22024.7505
$ws.Range("M136").Value = -4176.6489
$ws.Range("N136").Value = -27124.7505
